$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 41
$ws.Cells.Item(2, 12).Value = 'stimuli/img_aplao.png'
$ws.Cells.Item(2, 13).Value = 64.09090909090909
$ws.Cells.Item(2, 14).Value = 40.75757575757576
$ws.Cells.Item(2, 15).Value = 52.42424242424242
$ws.Cells.Item(2, 16).Value = 33
$ws.Cells.Item(2, 17).Value = 3
$ws.Cells.Item(2, 18).Value = 3
$ws.Cells.Item(2, 19).Value = 3
$ws.Cells.Item(2, 20).Value = 3
$ws.Cells.Item(2, 21).Value = 3
$ws.Cells.Item(2, 22).Value = 3

$ws.Cells.Item(3, 6).Value = 42

$ws.Cells.Item(4, 6).Value = 43

$ws.Cells.Item(5, 6).Value = 44

$ws.Cells.Item(6, 6).Value = 45
$ws.Cells.Item(6, 12).Value = 'stimuli/img_eatdk.png'
$ws.Cells.Item(6, 13).Value = 81.40625
$ws.Cells.Item(6, 14).Value = 61.375
$ws.Cells.Item(6, 15).Value = 71.390625
$ws.Cells.Item(6, 16).Value = 32
$ws.Cells.Item(6, 17).Value = 8
$ws.Cells.Item(6, 18).Value = 8
$ws.Cells.Item(6, 19).Value = 8
$ws.Cells.Item(6, 20).Value = 8
$ws.Cells.Item(6, 21).Value = 8
$ws.Cells.Item(6, 22).Value = 8

$ws.Cells.Item(7, 6).Value = 46
$ws.Cells.Item(7, 12).Value = 'stimuli/img_30vz5.png'
$ws.Cells.Item(7, 13).Value = 86.21212121212122
$ws.Cells.Item(7, 14).Value = 68.27272727272727
$ws.Cells.Item(7, 15).Value = 77.24242424242425
$ws.Cells.Item(7, 16).Value = 33
$ws.Cells.Item(7, 17).Value = 10
$ws.Cells.Item(7, 18).Value = 10
$ws.Cells.Item(7, 19).Value = 10
$ws.Cells.Item(7, 20).Value = 10
$ws.Cells.Item(7, 21).Value = 10
$ws.Cells.Item(7, 22).Value = 10

$ws.Cells.Item(8, 6).Value = 47

$ws.Cells.Item(9, 6).Value = 48
$ws.Cells.Item(9, 12).Value = 'stimuli/img_cv6mf.png'
$ws.Cells.Item(9, 13).Value = 66.8
$ws.Cells.Item(9, 14).Value = 42.08
$ws.Cells.Item(9, 15).Value = 54.44
$ws.Cells.Item(9, 16).Value = 25
$ws.Cells.Item(9, 17).Value = 4
$ws.Cells.Item(9, 18).Value = 4
$ws.Cells.Item(9, 19).Value = 4
$ws.Cells.Item(9, 20).Value = 4
$ws.Cells.Item(9, 21).Value = 4
$ws.Cells.Item(9, 22).Value = 4

$ws.Cells.Item(10, 6).Value = 49
$ws.Cells.Item(10, 12).Value = 'stimuli/img_p659z.png'
$ws.Cells.Item(10, 13).Value = 84.21621621621621
$ws.Cells.Item(10, 14).Value = 65.37837837837837
$ws.Cells.Item(10, 15).Value = 74.79729729729729
$ws.Cells.Item(10, 16).Value = 37
$ws.Cells.Item(10, 17).Value = 9
$ws.Cells.Item(10, 18).Value = 9
$ws.Cells.Item(10, 19).Value = 9
$ws.Cells.Item(10, 20).Value = 9
$ws.Cells.Item(10, 21).Value = 9
$ws.Cells.Item(10, 22).Value = 9

$ws.Cells.Item(11, 6).Value = 50

$ws.Cells.Item(12, 6).Value = 51
$ws.Cells.Item(12, 12).Value = 'stimuli/img_cnyac.png'
$ws.Cells.Item(12, 13).Value = 69.14705882352941
$ws.Cells.Item(12, 14).Value = 47.8235294117647
$ws.Cells.Item(12, 15).Value = 58.48529411764706
$ws.Cells.Item(12, 16).Value = 34
$ws.Cells.Item(12, 17).Value = 5
$ws.Cells.Item(12, 18).Value = 5
$ws.Cells.Item(12, 19).Value = 5
$ws.Cells.Item(12, 20).Value = 5
$ws.Cells.Item(12, 21).Value = 5
$ws.Cells.Item(12, 22).Value = 5

$ws.Cells.Item(13, 6).Value = 52

$ws.Cells.Item(14, 6).Value = 53

$ws.Cells.Item(15, 6).Value = 54

$ws.Cells.Item(16, 6).Value = 55
$ws.Cells.Item(16, 12).Value = 'stimuli/img_njmgp.png'
$ws.Cells.Item(16, 13).Value = 80.48148148148148
$ws.Cells.Item(16, 14).Value = 58.4074074074074
$ws.Cells.Item(16, 15).Value = 69.44444444444444
$ws.Cells.Item(16, 16).Value = 27
$ws.Cells.Item(16, 17).Value = 8
$ws.Cells.Item(16, 18).Value = 8
$ws.Cells.Item(16, 19).Value = 8
$ws.Cells.Item(16, 20).Value = 8
$ws.Cells.Item(16, 21).Value = 8
$ws.Cells.Item(16, 22).Value = 8

$ws.Cells.Item(17, 6).Value = 56
$ws.Cells.Item(17, 12).Value = 'stimuli/img_d8xbu.png'
$ws.Cells.Item(17, 13).Value = 91.36363636363636
$ws.Cells.Item(17, 14).Value = 73.18181818181819
$ws.Cells.Item(17, 15).Value = 82.27272727272728
$ws.Cells.Item(17, 16).Value = 33
$ws.Cells.Item(17, 17).Value = 10
$ws.Cells.Item(17, 18).Value = 10
$ws.Cells.Item(17, 19).Value = 10
$ws.Cells.Item(17, 20).Value = 10
$ws.Cells.Item(17, 21).Value = 10
$ws.Cells.Item(17, 22).Value = 10

$ws.Cells.Item(18, 6).Value = 57
$ws.Cells.Item(18, 12).Value = 'stimuli/img_p3hpc.png'
$ws.Cells.Item(18, 13).Value = 72.83333333333333
$ws.Cells.Item(18, 14).Value = 52.22222222222222
$ws.Cells.Item(18, 15).Value = 62.52777777777777
$ws.Cells.Item(18, 16).Value = 36
$ws.Cells.Item(18, 17).Value = 6
$ws.Cells.Item(18, 18).Value = 6
$ws.Cells.Item(18, 19).Value = 6
$ws.Cells.Item(18, 20).Value = 6
$ws.Cells.Item(18, 21).Value = 6
$ws.Cells.Item(18, 22).Value = 6

$ws.Cells.Item(19, 6).Value = 58
$ws.Cells.Item(19, 12).Value = 'stimuli/img_ce9vx.png'
$ws.Cells.Item(19, 13).Value = 75.90909090909091
$ws.Cells.Item(19, 14).Value = 57.12121212121212
$ws.Cells.Item(19, 15).Value = 66.51515151515152
$ws.Cells.Item(19, 17).Value = 7
$ws.Cells.Item(19, 18).Value = 7
$ws.Cells.Item(19, 19).Value = 7
$ws.Cells.Item(19, 20).Value = 7
$ws.Cells.Item(19, 21).Value = 7
$ws.Cells.Item(19, 22).Value = 7

$ws.Cells.Item(20, 6).Value = 59

$ws.Cells.Item(21, 6).Value = 60
$ws.Cells.Item(21, 12).Value = 'stimuli/img_t90e2.png'
$ws.Cells.Item(21, 13).Value = 83.0625
$ws.Cells.Item(21, 14).Value = 61.96875
$ws.Cells.Item(21, 15).Value = 72.515625
$ws.Cells.Item(21, 16).Value = 32
$ws.Cells.Item(21, 17).Value = 9
$ws.Cells.Item(21, 18).Value = 9
$ws.Cells.Item(21, 19).Value = 9
$ws.Cells.Item(21, 20).Value = 9
$ws.Cells.Item(21, 21).Value = 9
$ws.Cells.Item(21, 22).Value = 9

$ws.Cells.Item(22, 6).Value = 61

$ws.Cells.Item(23, 6).Value = 62

$ws.Cells.Item(24, 6).Value = 63
$ws.Cells.Item(24, 12).Value = 'stimuli/img_6nbgt.png'
$ws.Cells.Item(24, 13).Value = 78.45161290322581
$ws.Cells.Item(24, 14).Value = 57.83870967741935
$ws.Cells.Item(24, 15).Value = 68.14516129032258
$ws.Cells.Item(24, 16).Value = 31
$ws.Cells.Item(24, 17).Value = 7
$ws.Cells.Item(24, 18).Value = 7
$ws.Cells.Item(24, 19).Value = 7
$ws.Cells.Item(24, 20).Value = 7
$ws.Cells.Item(24, 21).Value = 7
$ws.Cells.Item(24, 22).Value = 7

$ws.Cells.Item(25, 6).Value = 64

$ws.Cells.Item(26, 6).Value = 65

$ws.Cells.Item(27, 6).Value = 66

$ws.Cells.Item(28, 6).Value = 67

$ws.Cells.Item(29, 6).Value = 68
$ws.Cells.Item(29, 12).Value = 'stimuli/img_es7o2.png'
$ws.Cells.Item(29, 13).Value = 52.48571428571429
$ws.Cells.Item(29, 14).Value = 27.54285714285714
$ws.Cells.Item(29, 15).Value = 40.01428571428572
$ws.Cells.Item(29, 16).Value = 35
$ws.Cells.Item(29, 17).Value = 2
$ws.Cells.Item(29, 18).Value = 2
$ws.Cells.Item(29, 19).Value = 2
$ws.Cells.Item(29, 20).Value = 2
$ws.Cells.Item(29, 21).Value = 2
$ws.Cells.Item(29, 22).Value = 2

$ws.Cells.Item(30, 6).Value = 69

$ws.Cells.Item(31, 6).Value = 70
$ws.Cells.Item(31, 12).Value = 'stimuli/img_yeh72.png'
$ws.Cells.Item(31, 13).Value = 68.66666666666667
$ws.Cells.Item(31, 14).Value = 45.21212121212121
$ws.Cells.Item(31, 15).Value = 56.93939393939394
$ws.Cells.Item(31, 16).Value = 33

$ws.Cells.Item(32, 6).Value = 71
$ws.Cells.Item(32, 12).Value = 'stimuli/img_7wul8.png'
$ws.Cells.Item(32, 13).Value = 43.03030303030303
$ws.Cells.Item(32, 14).Value = 25.54545454545455
$ws.Cells.Item(32, 15).Value = 34.28787878787879
$ws.Cells.Item(32, 16).Value = 33
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = 1
$ws.Cells.Item(32, 19).Value = 1
$ws.Cells.Item(32, 20).Value = 1
$ws.Cells.Item(32, 21).Value = 1
$ws.Cells.Item(32, 22).Value = 1

$ws.Cells.Item(33, 6).Value = 72

$ws.Cells.Item(34, 6).Value = 73
$ws.Cells.Item(34, 12).Value = 'stimuli/img_inqod.png'
$ws.Cells.Item(34, 13).Value = 70.84848484848484
$ws.Cells.Item(34, 14).Value = 50.63636363636363
$ws.Cells.Item(34, 15).Value = 60.74242424242424
$ws.Cells.Item(34, 17).Value = 5
$ws.Cells.Item(34, 18).Value = 5
$ws.Cells.Item(34, 19).Value = 5
$ws.Cells.Item(34, 20).Value = 5
$ws.Cells.Item(34, 21).Value = 5
$ws.Cells.Item(34, 22).Value = 5

$ws.Cells.Item(35, 6).Value = 74

$ws.Cells.Item(36, 6).Value = 75

$ws.Cells.Item(37, 6).Value = 76

$ws.Cells.Item(38, 6).Value = 77
$ws.Cells.Item(38, 12).Value = 'stimuli/img_iyxnj.png'
$ws.Cells.Item(38, 13).Value = 75.30555555555556
$ws.Cells.Item(38, 14).Value = 54.33333333333334
$ws.Cells.Item(38, 15).Value = 64.81944444444444
$ws.Cells.Item(38, 16).Value = 36
$ws.Cells.Item(38, 17).Value = 6
$ws.Cells.Item(38, 18).Value = 6
$ws.Cells.Item(38, 19).Value = 6
$ws.Cells.Item(38, 20).Value = 6
$ws.Cells.Item(38, 21).Value = 6
$ws.Cells.Item(38, 22).Value = 6

$ws.Cells.Item(39, 6).Value = 78
$ws.Cells.Item(39, 12).Value = 'stimuli/img_nyv2b.png'
$ws.Cells.Item(39, 13).Value = 11.91176470588235
$ws.Cells.Item(39, 14).Value = 6.852941176470588
$ws.Cells.Item(39, 15).Value = 9.382352941176471
$ws.Cells.Item(39, 16).Value = 34
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = 1
$ws.Cells.Item(39, 19).Value = 1
$ws.Cells.Item(39, 20).Value = 1
$ws.Cells.Item(39, 21).Value = 1
$ws.Cells.Item(39, 22).Value = 1

$ws.Cells.Item(40, 6).Value = 79

$ws.Cells.Item(41, 6).Value = 80
$ws.Cells.Item(41, 12).Value = 'stimuli/img_wyl6z.png'
$ws.Cells.Item(41, 13).Value = 59.8235294117647
$ws.Cells.Item(41, 14).Value = 36.23529411764706
$ws.Cells.Item(41, 15).Value = 48.02941176470588
$ws.Cells.Item(41, 16).Value = 34
$ws.Cells.Item(41, 17).Value = 3
$ws.Cells.Item(41, 18).Value = 3
$ws.Cells.Item(41, 19).Value = 3
$ws.Cells.Item(41, 20).Value = 3
$ws.Cells.Item(41, 21).Value = 3
$ws.Cells.Item(41, 22).Value = 3

$ws.Cells.Item(42, 6).Value = 81
$ws.Cells.Item(42, 12).Value = 'stimuli/img_ye5sl.png'
$ws.Cells.Item(42, 13).Value = 53.2258064516129
$ws.Cells.Item(42, 14).Value = 34.45161290322581
$ws.Cells.Item(42, 15).Value = 43.83870967741936
$ws.Cells.Item(42, 16).Value = 31
$ws.Cells.Item(42, 17).Value = 2
$ws.Cells.Item(42, 18).Value = 2
$ws.Cells.Item(42, 19).Value = 8
$ws.Cells.Item(42, 20).Value = 8
$ws.Cells.Item(42, 21).Value = 8
$ws.Cells.Item(42, 22).Value = 2
